# Update countries & provincias Spain
#
# 1) Three pairs/triples of country-name rows in column A were reordered
#    in the shared-string table of the source workbook, which (since each
#    worksheet row keeps pointing at a fixed shared-string slot) changes
#    which country name shows up on those particular rows:
#      - row 88/89   : El Salvador / Republica de Macedonia  -> swapped
#      - row 198/200 : Santa Lucia / Belice / Nueva Caledonia -> reversed
#                       (Belice, the middle one, stays put)
#      - row 215/216 : San Bartolome / Bonaire, San Eustaquio y Saba -> swapped
# 2) The "last updated" timestamp banner in A1 moved from 12:35 to 13:05.
# 3) A batch of per-country stat columns (B..H) were refreshed with newer
#    counts for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Country name reorders (column A) ---------------------------------
$ws.Range("A88").Value = "Republica de Macedonia"
$ws.Range("A89").Value = "El Salvador"

$ws.Range("A198").Value = "Nueva Caledonia"
$ws.Range("A199").Value = "Belice"
$ws.Range("A200").Value = "Santa Lucia"

$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"

# --- 2) Timestamp banner ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 13:05"

# --- 3) Refreshed per-country statistics (B..H) ----------------------------
$updates = @{
    14  = @{ B = 137724; C = 2023; D = 107713; E = 22560;           G = 34; H = 7451 }
    30  = @{ B = 30746;  C = 10;                E = 739;            G = 1;  H = 1907 }
    55  = @{                       D = 4400;    E = 4096 }
    83  = @{ B = 2406;   C = 5;    D = 1696;    E = 564;            G = 2;  H = 146 }
    88  = @{ B = 1999;   C = 21;   D = 1439;    E = 447;                    H = 113 }
    89  = @{ B = 1983;   C = 68;   D = 698;     E = 1250;                   H = 35 }
    99  = @{ B = 1469;   C = 1;    D = 1346;    E = 16 }
    126 = @{ B = 682;    C = 79;   D = 112;     E = 566;            G = 1;  H = 4 }
    130 = @{ B = 611;    C = 1;    D = 485;     E = 120 }
    134 = @{ B = 542;    C = 15;   D = 147;     E = 393 }
    146 = @{ B = 326;    C = 1;    D = 272;     E = 54 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
